$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.563.19"
$ws.Cells.Item(2, 5).Value = "  -0.59%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.628.51"
$ws.Cells.Item(3, 5).Value = "  -0.37%  "
$ws.Cells.Item(4, 5).Value = "  +0.21%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "211.83"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.521"
$ws.Cells.Item(6, 5).Value = "  -0.40%  "
$ws.Cells.Item(7, 5).Value = "  +0.20%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "22.95"
$ws.Cells.Item(8, 5).Value = "  -1.52%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.263"
$ws.Cells.Item(9, 5).Value = "  +0.71%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0611"
$ws.Cells.Item(10, 5).Value = "  +0.14%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0859"
$ws.Cells.Item(11, 5).Value = "  -3.41%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.861.95"
$ws.Cells.Item(12, 5).Value = "  -0.24%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.628.51"
$ws.Cells.Item(13, 5).Value = "  +0.01%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.03"
$ws.Cells.Item(14, 5).Value = "  -0.35%  "
$ws.Cells.Item(15, 5).Value = "  -1.09%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.89"
$ws.Cells.Item(16, 5).Value = "  +0.51%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "27.561.25"
$ws.Cells.Item(17, 5).Value = "  -0.54%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "228.11"
$ws.Cells.Item(18, 5).Value = "  -0.60%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.0₃0718"
$ws.Cells.Item(19, 5).Value = "  -0.48%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.52"
$ws.Cells.Item(20, 5).Value = "  -1.76%  "
$ws.Cells.Item(21, 5).Value = "  +0.14%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "10.63"
$ws.Cells.Item(22, 5).Value = "  +7.10%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.34"
$ws.Cells.Item(23, 5).Value = "  +1.07%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.12"
$ws.Cells.Item(24, 5).Value = "  +1.82%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "149.15"
$ws.Cells.Item(25, 5).Value = "  -1.32%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "6.86"
$ws.Cells.Item(26, 5).Value = "  -0.92%  "
$ws.Cells.Item(27, 5).Value = "  -0.96%  "
$ws.Cells.Item(28, 2).Value = "BinanceUSD"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.00"
$ws.Cells.Item(28, 5).Value = "  +0.22%  "
$ws.Cells.Item(29, 2).Value = "EthereumClassic"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.53"
$ws.Cells.Item(29, 5).Value = "  -0.41%  "
$ws.Cells.Item(30, 5).Value = "  -0.03%  "
$ws.Cells.Item(31, 5).Value = "  -0.57%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.27"
$ws.Cells.Item(32, 5).Value = "  -0.81%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.467.75"
$ws.Cells.Item(33, 5).Value = "  +0.18%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.07"
$ws.Cells.Item(34, 5).Value = "  -1.18%  "
$ws.Cells.Item(35, 5).Value = "  -1.35%  "
$ws.Cells.Item(36, 5).Value = "  -1.36%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.924"
$ws.Cells.Item(37, 5).Value = "  +0.17%  "
$ws.Cells.Item(38, 2).Value = "ImmutableX"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.558"
$ws.Cells.Item(38, 5).Value = "  -1.19%  "
$ws.Cells.Item(39, 2).Value = "ARBITRUM"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.874"
$ws.Cells.Item(39, 5).Value = "  -0.25%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0167"
$ws.Cells.Item(40, 5).Value = "  +0.05%  "
$ws.Cells.Item(41, 2).Value = "PaxDollar"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.00"
$ws.Cells.Item(41, 5).Value = "  +0.22%  "
$ws.Cells.Item(42, 2).Value = "WEMIXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.01"
$ws.Cells.Item(42, 5).Value = "  +0.41%  "
$ws.Cells.Item(43, 2).Value = "Aave"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "67.59"
$ws.Cells.Item(43, 5).Value = "  -1.28%  "
$ws.Cells.Item(44, 2).Value = "mCoin"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.47"
$ws.Cells.Item(44, 5).Value = "  +0.10%  "
$ws.Cells.Item(45, 2).Value = "MXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.21"
$ws.Cells.Item(45, 5).Value = "  -0.57%  "
$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "5.34"
$ws.Cells.Item(46, 5).Value = "  -3.34%  "
$ws.Cells.Item(47, 2).Value = "RocketPoolETH"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.770.20"
$ws.Cells.Item(47, 5).Value = "  -0.34%  "
$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.74"
$ws.Cells.Item(48, 5).Value = "  +2.76%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "87.26"
$ws.Cells.Item(49, 5).Value = "  +0.35%  "
$ws.Cells.Item(50, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0₆0106"
$ws.Cells.Item(50, 5).Value = "  +0.09%  "
$ws.Cells.Item(51, 2).Value = "Algorand"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0989"
$ws.Cells.Item(51, 5).Value = "  -0.07%  "
